$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Sema4d"
$ws.Range("C2").Value = "Erbb2"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 0.8730476666666668
$ws.Range("H2").Value = 2.619143
$ws.Range("I2").Value = 0.01740928848427011
$ws.Range("J2").Value = 0.01740928848427011
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 3.155977333333333
$ws.Range("N2").Value = 9.467932
$ws.Range("O2").Value = 0.3579027849973545
$ws.Range("P2").Value = 0.3579027849973545
$ws.Range("Q2").Value = 2.755318646919556
$ws.Range("R2").Value = 24.797867822276
$ws.Range("S2").Value = 0.006230832833342642
$ws.Range("T2").Value = 0.006230832833342644

$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Sema4d"
$ws.Range("C3").Value = "Erbb2"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 0.6666666666666666
$ws.Range("G3").Value = 0.8730476666666668
$ws.Range("H3").Value = 2.619143
$ws.Range("I3").Value = 0.01740928848427011
$ws.Range("J3").Value = 0.01740928848427011
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 3.165953666666667
$ws.Range("N3").Value = 9.497861
$ws.Range("O3").Value = 0.359034148472735
$ws.Range("P3").Value = 0.359034148472735
$ws.Range("Q3").Value = 2.764028461458111
$ws.Range("R3").Value = 24.876256153123
$ws.Range("S3").Value = 0.006250529066466108
$ws.Range("T3").Value = 0.006250529066466109

$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Sema4d"
$ws.Range("C4").Value = "Erbb2"
$ws.Range("D4").Value = "sCs"
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 0.6666666666666666
$ws.Range("G4").Value = 0.8730476666666668
$ws.Range("H4").Value = 2.619143
$ws.Range("I4").Value = 0.01740928848427011
$ws.Range("J4").Value = 0.01740928848427011
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 2.496042666666666
$ws.Range("N4").Value = 7.488128
$ws.Range("O4").Value = 0.2830630665299106
$ws.Range("P4").Value = 0.2830630665299106
$ws.Range("Q4").Value = 2.179164226033778
$ws.Range("R4").Value = 19.612478034304
$ws.Range("S4").Value = 0.004927926584461356
$ws.Range("T4").Value = 0.004927926584461357

$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Sema4d"
$ws.Range("C5").Value = "Erbb2"
$ws.Range("D5").Value = "ECs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 1.252512666666667
$ws.Range("H5").Value = 3.757538
$ws.Range("I5").Value = 0.02497613266347325
$ws.Range("J5").Value = 0.02497613266347325
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 3.155977333333333
$ws.Range("N5").Value = 9.467932
$ws.Range("O5").Value = 0.3579027849973545
$ws.Range("P5").Value = 0.3579027849973545
$ws.Range("Q5").Value = 3.952901585712889
$ws.Range("R5").Value = 35.576114271416
$ws.Range("S5").Value = 0.00893902743872047
$ws.Range("T5").Value = 0.00893902743872047

$ws.Range("A6").Value = "FAPs"
$ws.Range("B6").Value = "Sema4d"
$ws.Range("C6").Value = "Erbb2"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 1.252512666666667
$ws.Range("H6").Value = 3.757538
$ws.Range("I6").Value = 0.02497613266347325
$ws.Range("J6").Value = 0.02497613266347325
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 3.165953666666667
$ws.Range("N6").Value = 9.497861
$ws.Range("O6").Value = 0.359034148472735
$ws.Range("P6").Value = 0.359034148472735
$ws.Range("Q6").Value = 3.965397069579778
$ws.Range("R6").Value = 35.688573626218
$ws.Range("S6").Value = 0.008967284522972182
$ws.Range("T6").Value = 0.008967284522972182

$ws.Range("A7").Value = "FAPs"
$ws.Range("B7").Value = "Sema4d"
$ws.Range("C7").Value = "Erbb2"
$ws.Range("D7").Value = "sCs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 1.252512666666667
$ws.Range("H7").Value = 3.757538
$ws.Range("I7").Value = 0.02497613266347325
$ws.Range("J7").Value = 0.02497613266347325
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 2.496042666666666
$ws.Range("N7").Value = 7.488128
$ws.Range("O7").Value = 0.2830630665299106
$ws.Range("P7").Value = 0.2830630665299106
$ws.Range("Q7").Value = 3.126325056540444
$ws.Range("R7").Value = 28.136925508864
$ws.Range("S7").Value = 0.007069820701780604
$ws.Range("T7").Value = 0.007069820701780604

$ws.Range("A8").Value = "M2"
$ws.Range("B8").Value = "Sema4d"
$ws.Range("C8").Value = "Erbb2"
$ws.Range("D8").Value = "ECs"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 45.633473
$ws.Range("H8").Value = 136.900419
$ws.Range("I8").Value = 0.9099689814525027
$ws.Range("J8").Value = 0.9099689814525027
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 3.155977333333333
$ws.Range("N8").Value = 9.467932
$ws.Range("O8").Value = 0.3579027849973545
$ws.Range("P8").Value = 0.3579027849973545
$ws.Range("Q8").Value = 144.0182064292786
$ws.Range("R8").Value = 1296.163857863508
$ws.Range("S8").Value = 0.3256804327230567
$ws.Range("T8").Value = 0.3256804327230567

$ws.Range("A9").Value = "M2"
$ws.Range("B9").Value = "Sema4d"
$ws.Range("C9").Value = "Erbb2"
$ws.Range("D9").Value = "FAPs"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 45.633473
$ws.Range("H9").Value = 136.900419
$ws.Range("I9").Value = 0.9099689814525027
$ws.Range("J9").Value = 0.9099689814525027
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 3.165953666666667
$ws.Range("N9").Value = 9.497861
$ws.Range("O9").Value = 0.359034148472735
$ws.Range("P9").Value = 0.359034148472735
$ws.Range("Q9").Value = 144.4734611670843
$ws.Range("R9").Value = 1300.261150503759
$ws.Range("S9").Value = 0.3267099383924013
$ws.Range("T9").Value = 0.3267099383924013

$ws.Range("A10").Value = "M2"
$ws.Range("B10").Value = "Sema4d"
$ws.Range("C10").Value = "Erbb2"
$ws.Range("D10").Value = "sCs"
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 45.633473
$ws.Range("H10").Value = 136.900419
$ws.Range("I10").Value = 0.9099689814525027
$ws.Range("J10").Value = 0.9099689814525027
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 2.496042666666666
$ws.Range("N10").Value = 7.488128
$ws.Range("O10").Value = 0.2830630665299106
$ws.Range("P10").Value = 0.2830630665299106
$ws.Range("Q10").Value = 113.9030956361813
$ws.Range("R10").Value = 1025.127860725632
$ws.Range("S10").Value = 0.2575786103370448
$ws.Range("T10").Value = 0.2575786103370448

$ws.Range("A11").Value = "sCs"
$ws.Range("B11").Value = "Sema4d"
$ws.Range("C11").Value = "Erbb2"
$ws.Range("D11").Value = "ECs"
$ws.Range("E11").Value = 3
$ws.Range("F11").Value = 1
$ws.Range("G11").Value = 2.389349666666667
$ws.Range("H11").Value = 7.168049000000001
$ws.Range("I11").Value = 0.04764559739975398
$ws.Range("J11").Value = 0.04764559739975399
$ws.Range("K11").Value = 3
$ws.Range("L11").Value = 1
$ws.Range("M11").Value = 3.155977333333333
$ws.Range("N11").Value = 9.467932
$ws.Range("O11").Value = 0.3579027849973545
$ws.Range("P11").Value = 0.3579027849973545
$ws.Range("Q11").Value = 7.540733389407555
$ws.Range("R11").Value = 67.866600504668
$ws.Range("S11").Value = 0.01705249200223466
$ws.Range("T11").Value = 0.01705249200223466

$ws.Range("A12").Value = "sCs"
$ws.Range("B12").Value = "Sema4d"
$ws.Range("C12").Value = "Erbb2"
$ws.Range("D12").Value = "FAPs"
$ws.Range("E12").Value = 3
$ws.Range("F12").Value = 1
$ws.Range("G12").Value = 2.389349666666667
$ws.Range("H12").Value = 7.168049000000001
$ws.Range("I12").Value = 0.04764559739975398
$ws.Range("J12").Value = 0.04764559739975399
$ws.Range("K12").Value = 3
$ws.Range("L12").Value = 1
$ws.Range("M12").Value = 3.165953666666667
$ws.Range("N12").Value = 9.497861
$ws.Range("O12").Value = 0.359034148472735
$ws.Range("P12").Value = 0.359034148472735
$ws.Range("Q12").Value = 7.564570338132112
$ws.Range("R12").Value = 68.08113304318901
$ws.Range("S12").Value = 0.01710639649089543
$ws.Range("T12").Value = 0.01710639649089543

$ws.Range("A13").Value = "sCs"
$ws.Range("B13").Value = "Sema4d"
$ws.Range("C13").Value = "Erbb2"
$ws.Range("D13").Value = "sCs"
$ws.Range("E13").Value = 3
$ws.Range("F13").Value = 1
$ws.Range("G13").Value = 2.389349666666667
$ws.Range("H13").Value = 7.168049000000001
$ws.Range("I13").Value = 0.04764559739975398
$ws.Range("J13").Value = 0.04764559739975399
$ws.Range("K13").Value = 3
$ws.Range("L13").Value = 1
$ws.Range("M13").Value = 2.496042666666666
$ws.Range("N13").Value = 7.488128
$ws.Range("O13").Value = 0.2830630665299106
$ws.Range("P13").Value = 0.2830630665299106
$ws.Range("Q13").Value = 5.963918713585778
$ws.Range("R13").Value = 53.675268422272
$ws.Range("S13").Value = 0.0134867089066239
$ws.Range("T13").Value = 0.0134867089066239
